$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15.66

$ws.Range("B3").Value = 21.54
$ws.Range("C3").Value = 16.37

$ws.Range("C4").Value = 16.03

$ws.Range("C5").Value = 15.9

$ws.Range("C7").Value = 15.34

$ws.Range("C9").Value = 13.05

$ws.Range("C10").Value = 12.84

$ws.Range("C11").Value = 13.89

$ws.Range("C13").Value = 14.59

$ws.Range("C15").Value = 12.89

$ws.Range("C17").Value = 13.84

$ws.Range("C18").Value = 14.48

$ws.Range("C21").Value = 13.64

$ws.Range("C22").Value = 18.31

$ws.Range("C24").Value = 22.68
